$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 105; everything from row 105 down shifts to 107 down.
$ws.Rows("105:106").Insert()

# Shared metadata columns (A-C, E-J) are identical across all data rows in this sheet.
$marketId = 2
$market   = "Comercializadora del Agro de Limarí"
$region   = "Coquimbo"
$codreg   = 4
$tipo     = "Fruta"
$prodId   = 100103
$producto = "Frutos de hueso (carozo)"
$catId    = 100103004
$categoria = "Durazno"

# New row 105: September Snow / Primera, bins lot dated 2022-03-16 (serial 44636)
$ws.Cells.Item(105, 1).Value = $marketId
$ws.Cells.Item(105, 2).Value = $market
$ws.Cells.Item(105, 3).Value = $region
$ws.Cells.Item(105, 4).Value = 44636
$ws.Cells.Item(105, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(105, 5).Value = $codreg
$ws.Cells.Item(105, 6).Value = $tipo
$ws.Cells.Item(105, 7).Value = $prodId
$ws.Cells.Item(105, 8).Value = $producto
$ws.Cells.Item(105, 9).Value = $catId
$ws.Cells.Item(105, 10).Value = $categoria
$ws.Cells.Item(105, 11).Value = "September Snow"
$ws.Cells.Item(105, 12).Value = "Primera"
$ws.Cells.Item(105, 13).Value = 20
$ws.Cells.Item(105, 14).Value = 355000
$ws.Cells.Item(105, 15).Value = 360000
$ws.Cells.Item(105, 16).Value = 357500
$ws.Cells.Item(105, 17).Value = "$/bins (400 kilos)"
$ws.Cells.Item(105, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(105, 19).Value = 894
$ws.Cells.Item(105, 20).Value = 400

# New row 106: September Snow / Segunda, same bins lot
$ws.Cells.Item(106, 1).Value = $marketId
$ws.Cells.Item(106, 2).Value = $market
$ws.Cells.Item(106, 3).Value = $region
$ws.Cells.Item(106, 4).Value = 44636
$ws.Cells.Item(106, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(106, 5).Value = $codreg
$ws.Cells.Item(106, 6).Value = $tipo
$ws.Cells.Item(106, 7).Value = $prodId
$ws.Cells.Item(106, 8).Value = $producto
$ws.Cells.Item(106, 9).Value = $catId
$ws.Cells.Item(106, 10).Value = $categoria
$ws.Cells.Item(106, 11).Value = "September Snow"
$ws.Cells.Item(106, 12).Value = "Segunda"
$ws.Cells.Item(106, 13).Value = 20
$ws.Cells.Item(106, 14).Value = 305000
$ws.Cells.Item(106, 15).Value = 310000
$ws.Cells.Item(106, 16).Value = 307500
$ws.Cells.Item(106, 17).Value = "$/bins (400 kilos)"
$ws.Cells.Item(106, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(106, 19).Value = 769
$ws.Cells.Item(106, 20).Value = 400
